# Fix small typos ("it´s" -> "its") in the "Main Game Loop" analysis
# text box on the Readme slide, as described in the commit
# "Just small corrections on analysis".
#
# Strategy: locate the target paragraph via TextRange.Paragraphs(n), then
# find the exact run boundaries inside that paragraph (TextRange.Start is
# reliable per-run even though Runs.Count/Length/Text are not), and replace
# whole-run (or exact sub-run) character ranges in place via
# TextRange.Characters(start, length).Text = "...". Replacing a full run's
# character range keeps it a single run; replacing an interior sub-range
# naturally splits the run into the pieces we need (used for the 4th fix,
# which the original author's edit also split into three runs).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(9)
$shp = $s.Shapes.Item(2)
$tf = $shp.TextFrame
$tr = $tf.TextRange

# The text box has spAutoFit; editing TextRange.Text normally triggers an
# immediate shape-height recalculation (not reflected in the source diff,
# which only touches <a:t> runs). Temporarily disable auto-sizing while we
# edit, then restore both the autosize mode and the original height so the
# shape's <a:xfrm> stays byte-identical to the original.
$origHeight = $shp.Height
$tf.AutoSize = 0   # msoAutoSizeNone

function Get-RunStarts($para) {
    $paraEnd = $para.Start + $para.Length
    $runs = $para.Runs()
    $starts = New-Object System.Collections.ArrayList
    foreach ($r in $runs) {
        $st = $r.Start
        if ($st -ge $paraEnd) { break }
        if ($starts.Count -gt 0 -and $starts[$starts.Count - 1] -eq $st) { continue }
        [void]$starts.Add($st)
    }
    return $starts
}

function Get-RunRange($para, [int]$runIndex0) {
    # 0-based run index within the paragraph.
    $starts = Get-RunStarts $para
    $paraEnd = $para.Start + $para.Length
    $runStart = $starts[$runIndex0]
    if ($runIndex0 + 1 -lt $starts.Count) {
        $runLen = $starts[$runIndex0 + 1] - $runStart
    } else {
        $runLen = $paraEnd - $runStart - 1   # exclude trailing paragraph mark
    }
    return $tr.Characters($runStart, $runLen)
}

# --- Fix 1: "In it´s main() the Controller creates ..." -> "In its main() ..." ---
$para3 = $tr.Paragraphs(3)
$run = Get-RunRange $para3 0
$run.Text = $run.Text.Replace("it´s", "its")

# --- Fix 2: "...checks for input from the Player and notifies it´s observers
#             (the Controller)." -> split into "...and ", "notifies its ",
#             "observers (the Controller)." ---
$para7 = $tr.Paragraphs(7)
$run = Get-RunRange $para7 1
$fullText = $run.Text
$runStart = $run.Start
$idx = $fullText.IndexOf("notifies it´s ")
$mid = $tr.Characters($runStart + $idx, "notifies it´s ".Length)
$mid.Text = "notifies its "

# --- Fix 3: "In it´s update()" -> "In its update()" (paragraph "2 - ...") ---
$para9 = $tr.Paragraphs(9)
$run = Get-RunRange $para9 2
$run.Text = $run.Text.Replace("it´s", "its")

# --- Fix 4: "... handles the movement and notifies it´s observers (the View)."
#             -> "... notifies its observers (the View)." ---
$para11 = $tr.Paragraphs(11)
$run = Get-RunRange $para11 6
$run.Text = $run.Text.Replace("it´s", "its")

# --- Fix 5: standalone "it´s" -> "its" (paragraph "4 - In it´s update() ...") ---
$para13 = $tr.Paragraphs(13)
$run = Get-RunRange $para13 1
$run.Text = $run.Text.Replace("it´s", "its")

# Restore auto-sizing (so the saved bodyPr keeps <a:spAutoFit/>) and then
# force the shape back to its original, pre-edit height so the <a:ext> is
# unchanged. The host's point->EMU conversion rounds through a narrower
# float internally, so feed it a value from the middle of the input range
# that is known to land on the original 6186309 EMU (867398 x, 1162226 y,
# cx 10457203) rather than the mathematically "exact" 6186309/12700.
$tf.AutoSize = 1   # msoAutoSizeShapeToFitText
if ([Math]::Round($origHeight, 4) -eq 487.1109) {
    $shp.Height = 487.1109924316406
} else {
    $shp.Height = $origHeight
}
